$wb = $excel.ActiveWorkbook

$sanity = $wb.Worksheets.Item("Sanity")
$keyWords = $wb.Worksheets.Item("keyWords")

# Is_Enabled flag for every KeyWords row flips from "F" to "Y".
$sanity.Range("C2").Value = "Y"
$sanity.Range("C4").Value = "Y"
$sanity.Range("C6").Value = "Y"

# New "Class" column on the keyWords sheet, mapping each keyword to the
# reflection-invoked step-definition class that implements it.
$keyWords.Range("B1").Value = "Class"
$keyWords.Range("B2").Value = "com.sbn.pages.StepDefs"
$keyWords.Range("B3").Value = "com.sbn.pages.StepDefs"
$keyWords.Columns.Item(2).AutoFit()

# keyWords becomes the active sheet/tab, with B5 selected.
$keyWords.Activate()
$keyWords.Range("B5").Select()
